# Add a new first column ("testcase_number") to the User_Credentials sheet,
# shifting the existing "id"/"password" data one column to the right, and
# populate the new column with test-case ids (tc_001..tc_004).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before column A. This shifts the existing columns
#    A->B and B->C, carrying cell values *and* styles with them.
$ws.Columns("A").Insert()

# 2) Set the column widths: new A & B both take the old column-A width,
#    and C takes the old column-B width.
$ws.Columns("A:B").ColumnWidth = 23.73
$ws.Columns("C").ColumnWidth = 30.05

# 3) Fill in the new "testcase_number" column.
$ws.Range("A1").Value = "testcase_number"
$ws.Range("A2").Value = "tc_001"
$ws.Range("A3").Value = "tc_002"
$ws.Range("A4").Value = "tc_003"
$ws.Range("A5").Value = "tc_004"

# 4) The column insert does not move the worksheet's Hyperlinks along with
#    the cells they were anchored to, so the hyperlink anchors are still
#    pointing at column B. Re-create them on column C (where the
#    "balaPB@1234" / "Bala@130105" display text now lives), then restore
#    the plain "Normal" cell style so re-adding the hyperlink doesn't leave
#    the cells with Excel's blue/underlined Hyperlink formatting.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:balaPB@1234", "", "", "balaPB@1234")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:balaPB@1234", "", "", "balaPB@1234")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:Bala@130105", "", "", "Bala@130105")

$ws.Range("C2").Style = "Normal"
$ws.Range("C4").Style = "Normal"
$ws.Range("C5").Style = "Normal"

# 5) Put the selection where the diff shows it ended up.
$ws.Range("A5").Select()
